$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.674.61"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "2.173.24"
$ws.Range("E3").Value = "  -2.75%  "
$ws.Range("E4").Value = "  -0.06%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "238.56"
$c.ClearFormats()
$ws.Range("E5").Value = "  -1.81%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.612"
$c.ClearFormats()
$ws.Range("E6").Value = "  -2.25%  "
$ws.Range("E7").Value = "  -1.91%  "
$ws.Range("E8").Value = "  -0.11%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.583"
$c.ClearFormats()
$ws.Range("E9").Value = "  -2.87%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "40.24"
$c.ClearFormats()
$ws.Range("E10").Value = "  -4.78%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0911"
$c.ClearFormats()
$ws.Range("E11").Value = "  -4.79%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "54.51"
$c.ClearFormats()
$ws.Range("E12").Value = "  -3.62%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "6.73"
$c.ClearFormats()
$ws.Range("E13").Value = "  -2.90%  "
$ws.Range("E14").Value = "  -3.11%  "
$ws.Range("D15").Value = "2.500.59"
$ws.Range("E15").Value = "  -2.71%  "
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").Value = "2.186.87"
$ws.Range("E17").Value = "  -2.05%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.783"
$c.ClearFormats()
$ws.Range("E18").Value = "  -6.53%  "
$ws.Range("D19").Value = "41.647.69"
$ws.Range("E19").Value = "  -1.03%  "
$ws.Range("E20").Value = "  -2.00%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.81"
$c.ClearFormats()
$ws.Range("E22").Value = "  -6.56%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "10.10"
$c.ClearFormats()
$ws.Range("E23").Value = "  -11.89%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "226.45"
$c.ClearFormats()
$ws.Range("E24").Value = "  -1.42%  "
$ws.Range("E25").Value = "  +0.99%  "
$ws.Range("E26").Value = "  +0.14%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "10.77"
$c.ClearFormats()
$ws.Range("E27").Value = "  -5.30%  "
$ws.Range("E28").Value = "  -9.92%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.20"
$c.ClearFormats()
$ws.Range("E29").Value = "  -3.34%  "
$ws.Range("E30").Value = "  -1.08%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "170.72"
$c.ClearFormats()
$ws.Range("E31").Value = "  +1.99%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "19.88"
$c.ClearFormats()
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "32.86"
$c.ClearFormats()
$ws.Range("E33").Value = "  +10.17%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.0777"
$c.ClearFormats()
$ws.Range("E34").Value = "  -3.40%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "5.31"
$c.ClearFormats()
$ws.Range("E35").Value = "  -5.51%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.120"
$c.ClearFormats()
$ws.Range("E36").Value = "  -3.48%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "4.31"
$c.ClearFormats()
$ws.Range("E37").Value = "  -1.33%  "
$ws.Range("E38").Value = "  -6.30%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.0310"
$c.ClearFormats()
$ws.Range("E39").Value = "  +2.07%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "12.11"
$c.ClearFormats()
$ws.Range("E40").Value = "  -8.36%  "
$ws.Range("E41").Value = "  -1.71%  "
$ws.Range("E42").Value = "  -5.72%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "59.41"
$c.ClearFormats()
$ws.Range("E43").Value = "  -8.03%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.190"
$c.ClearFormats()
$ws.Range("E44").Value = "  -4.84%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "8.47"
$c.ClearFormats()
$ws.Range("E45").Value = "  -2.82%  "
$ws.Range("E46").Value = "  -3.48%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "97.39"
$c.ClearFormats()
$ws.Range("E47").Value = "  -6.78%  "
$ws.Range("E48").Value = "  -4.58%  "
$ws.Range("E49").Value = "  -4.94%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "2.20"
$c.ClearFormats()
$ws.Range("E50").Value = "  -6.04%  "
$ws.Range("E51").Value = "  -2.03%  "
